$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.240.49'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.04%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.911.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.85%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.86'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.26%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.82'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.26%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.61%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.63%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.35%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000344'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.17%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.05'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.87%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.94'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.96%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.537.31'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.84%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.016.45'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.95%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.25'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.35%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.24%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.92'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.34%  '

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.89%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.359.54'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.99%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '433.53'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.35%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.20%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.87'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.60%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.07%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.33'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +14.11%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.27'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +10.12%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.04%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.07'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.59%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.72'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.95%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '720.62'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.30%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.75'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.34%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.06%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.25'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +16.42%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '41.56'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.91%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0869'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '60.34'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.92%  '

# Row 38
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.148'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.32%  '

# Row 39
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.398'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +18.23%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.13%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0483'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.28%  '

# Row 42
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.92'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +14.80%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.14'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.75%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.93%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.47%  '

# Row 46
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.34'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.91%  '

# Row 47
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.07%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0354'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +38.84%  '

# Row 49
$ws.Range("B49").Value = 'LidoDAOToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.39'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.77%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.13'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.49%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '144.59'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.53%  '
